# Lengths.xlsx update
#
# - Insert a new column A ("Component" grouping) before the existing data;
#   this shifts the old Part/Height/Contribution/Notes columns one column
#   to the right (A->B, B->C, C->D, D->E) along with their formulas.
# - Rename a handful of "Part" labels (now living in column B) to be more
#   descriptive.
# - Old header "Component" (now in B2) becomes "Part"; new A2 header
#   becomes "Component".
# - Add grouping labels in the new column A for the first row of every
#   component group (Coupler / S2 / S1 / S1-Elect / Elect / Nose Cone).
# - Add a new warning note row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new, blank column before column A. Excel shifts every
# existing column (and the formulas/relative refs inside it) one slot to
# the right.
$ws.Columns("A:A").Insert()

# --- Rename the inter-mount related parts (column B, after the shift) ---
$ws.Range("B6").Value = "s2 inter-mount Disc"
$ws.Range("B7").Value = "s1 inter-mount Disc"
$ws.Range("B5").Value = "Tube inter-mount"
$ws.Range("B8").Value = "Tube inter-mount"

# --- Header row (row 2) ---
$ws.Range("B2").Value = "Part"
$ws.Range("A2").Value = "Component"
$ws.Range("A2").Font.Bold = $true

# --- Other renamed Part label ---
$ws.Range("B3").Value = "Coupler Tube"

# --- New column A grouping labels (one per component group) ---
$ws.Range("A4").Value = "S2 "
$ws.Range("A7").Value = "S1"
$ws.Range("A10").Value = "S1-Elect"
$ws.Range("A13").Value = "Elect"
$ws.Range("A3").Value = "Coupler"
$ws.Range("A16").Value = "Nose Cone"

# Size column A to fit its (short) contents, like the other bestFit columns.
$ws.Columns("A:A").AutoFit()

# --- New warning note row ---
$ws.Range("A20").Value = "THESE VALUES ARE GOING TO CHANGE!"

# --- Selection matches the saved workbook state ---
$ws.Range("A21").Select()

$wb.Save()
